# Apply the "balanços concatenados" edit: the O column (31/12/2009) is a
# duplicate period being removed (cleared, matching the already-blank D
# column / 31/12/2006), a handful of previously-zero-filled trailing
# columns (S:X) on subtotal/placeholder rows are cleared too, and a few
# cells get small recalculated numeric corrections.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells that become blank (matches the source's "empty inlineStr" cells,
#     i.e. a cleared cell with no value) ---------------------------------
$rangesToClear = @(
    "O57", "S57:X57",
    "O58", "S58:X58",
    "O59",
    "O60",
    "O61",
    "O62",
    "O63",
    "B64:C64", "E64:S64",
    "O65",
    "O66",
    "O67",
    "O68",
    "O69",
    "O70",
    "O71", "S71:X71",
    "O72", "S72:X72",
    "O73", "S73:X73",
    "O74",
    "O75",
    "O76",
    "O77", "S77:X77",
    "O78", "S78:X78",
    "B79:C79", "E79:O79", "Q79",
    "O80"
)

foreach ($rng in $rangesToClear) {
    $ws.Range($rng).ClearContents()
}

# --- Small numeric corrections (recalculated values) ---------------------
$ws.Range("S59").Value = 1099186.944
$ws.Range("W59").Value = 697025.9840000001
$ws.Range("S61").Value = 224399.024
$ws.Range("H69").Value = -45220
$ws.Range("W74").Value = 153906
$ws.Range("S79").Value = -102268
$ws.Range("W79").Value = -96713
